# daily auto push: 2026-02-23 19:19 UTC
# Insert two missing time-slot rows (2026/02/23 23:00 and 2026/02/24 02:00)
# right after the existing 2026/02/23 rows, shifting the remainder of the
# table (2026/12/29 .. 2027/01/05) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 843..884 down to 845..886, creating two empty rows at 843:844
$ws.Rows("843:844").Insert()

# Keep the date / weekday columns as plain text (they are stored as text
# in this workbook, not as real dates) so Excel does not auto-convert the
# "yyyy/mm/dd" strings into date serial numbers.
$ws.Range("A843:B844").NumberFormat = "@"

$ws.Range("A843").Value = "2026/02/23"
$ws.Range("B843").Value = "月"
$ws.Range("C843").Value = 23
$ws.Range("D843").Value = 29

$ws.Range("A844").Value = "2026/02/24"
$ws.Range("B844").Value = "火"
$ws.Range("C844").Value = 2
$ws.Range("D844").Value = 34
